$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the image_url entries in column E one row down starting at row 17:
# E17 currently holds the boomplaymusic url, E18 holds the pixabay url.
# After the edit, E17 should hold the pixabay url, E18 should be empty,
# and E19 should hold the boomplaymusic url.

$ws.Range("E17").Value = $ws.Range("E18").Value2
$ws.Range("E18").ClearContents()
$ws.Range("E19").Value = "https://source.boomplaymusic.com/group10/M00/06/08/9b2ee4c30c40406b981a55ee89a00ad1_464_464.webp"

# Update the active selection to match the new state
$ws.Range("D19").Select()
